# Apply updated market-board price figures to the leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1211692.2
$ws.Range("J17").Value = 1211692.2
$ws.Range("L17").Value = 3635076.6
$ws.Range("N17").Value = -3635412.6
$ws.Range("H33").Value = 390.82144
$ws.Range("I33").Value = 175.09091
$ws.Range("J33").Value = 1181.8334
$ws.Range("K33").Value = 175.09091
$ws.Range("L33").Value = 1181.8334
$ws.Range("M33").Value = 53.90908999999999
$ws.Range("N33").Value = -1639.8334
$ws.Range("H96").Value = 269.20834
$ws.Range("I96").Value = 231.16667
$ws.Range("J96").Value = 383.33334
$ws.Range("K96").Value = 693.50001
$ws.Range("L96").Value = 1150.00002
$ws.Range("M96").Value = 679.49999
$ws.Range("N96").Value = -3896.00002
$ws.Range("H116").Value = 6324.488
$ws.Range("I116").Value = 9300.333000000001
$ws.Range("J116").Value = 4607.654
$ws.Range("K116").Value = 9300.333000000001
$ws.Range("L116").Value = 4607.654
$ws.Range("M116").Value = -5858.333000000001
$ws.Range("N116").Value = -11491.654
$ws.Range("H129").Value = 956.5333000000001
$ws.Range("I129").Value = 698
$ws.Range("K129").Value = 2094
$ws.Range("M129").Value = 2906
$ws.Range("H132").Value = 130146.54
$ws.Range("I132").Value = 1721.1846
$ws.Range("J132").Value = 772273.3
$ws.Range("K132").Value = 5163.5538
$ws.Range("L132").Value = 2316819.9
$ws.Range("M132").Value = -2633.5538
$ws.Range("N132").Value = -2321879.9
$ws.Range("H134").Value = 47636
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 47636
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 47636
$ws.Range("N134").Value = -57776
$ws.Range("M134").ClearContents()
$ws.Range("H137").Value = 928.7838
$ws.Range("I137").Value = 590.875
$ws.Range("J137").Value = 1552.6154
$ws.Range("K137").Value = 1772.625
$ws.Range("L137").Value = 4657.8462
$ws.Range("M137").Value = 777.375
$ws.Range("N137").Value = -9757.8462
$ws.Range("H140").Value = 60738.77
$ws.Range("I140").Value = 60709
$ws.Range("J140").Value = 60741.25
$ws.Range("K140").Value = 60709
$ws.Range("L140").Value = 60741.25
$ws.Range("M140").Value = -55529
$ws.Range("N140").Value = -71101.25
$ws.Range("H141").Value = 2841
$ws.Range("I141").Value = 2525
$ws.Range("K141").Value = 7575
$ws.Range("M141").Value = -2395

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 143453.7
$ws.Range("I132").Value = 22857.746
$ws.Range("J132").Value = 419819.4
$ws.Range("K132").Value = 68573.238
$ws.Range("L132").Value = 1259458.2
$ws.Range("M132").Value = -66043.238
$ws.Range("N132").Value = -1264518.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 532054.6
$ws.Range("I105").Value = 885457.7
$ws.Range("K105").Value = 885457.7
$ws.Range("M105").Value = -883710.7
$ws.Range("H134").Value = 62565400
$ws.Range("I134").Value = 3035.0908
$ws.Range("J134").Value = 200202610
$ws.Range("K134").Value = 9105.2724
$ws.Range("L134").Value = 600607830
$ws.Range("M134").Value = -6570.2724
$ws.Range("N134").Value = -600612900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 230.25
$ws.Range("I23").Value = 57.75
$ws.Range("J23").Value = 316.5
$ws.Range("K23").Value = 173.25
$ws.Range("L23").Value = 949.5
$ws.Range("M23").Value = 61.75
$ws.Range("N23").Value = -1419.5
$ws.Range("H113").Value = 511.57144
$ws.Range("I113").Value = 510.4
$ws.Range("J113").Value = 514.5
$ws.Range("K113").Value = 1531.2
$ws.Range("L113").Value = 1543.5
$ws.Range("M113").Value = 638.8000000000002
$ws.Range("N113").Value = -5883.5
$ws.Range("H131").Value = 286484.28
$ws.Range("J131").Value = 303809.1
$ws.Range("L131").Value = 911427.2999999999
$ws.Range("N131").Value = -921507.2999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1691.7693
$ws.Range("I126").Value = 1455
$ws.Range("J126").Value = 1856.4783
$ws.Range("K126").Value = 4365
$ws.Range("L126").Value = 5569.4349
$ws.Range("M126").Value = -1895
$ws.Range("N126").Value = -10509.4349
$ws.Range("H141").Value = 43253.363
$ws.Range("J141").Value = 43253.363
$ws.Range("L141").Value = 43253.363
$ws.Range("N141").Value = -53613.363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2546.5833
$ws.Range("I7").Value = 2545.4644
$ws.Range("J7").Value = 2550.5
$ws.Range("K7").Value = 2545.4644
$ws.Range("L7").Value = 2550.5
$ws.Range("M7").Value = -2433.4644
$ws.Range("N7").Value = -2774.5
$ws.Range("H40").Value = 39820.555
$ws.Range("I40").Value = 2002.5
$ws.Range("J40").Value = 55743.95
$ws.Range("K40").Value = 2002.5
$ws.Range("L40").Value = 55743.95
$ws.Range("M40").Value = -1866.5
$ws.Range("N40").Value = -56015.95
$ws.Range("H122").Value = 33335830
$ws.Range("I122").Value = 40002530
$ws.Range("J122").Value = 22224664
$ws.Range("K122").Value = 120007590
$ws.Range("L122").Value = 66673992
$ws.Range("M122").Value = -120005140
$ws.Range("N122").Value = -66678892
$ws.Range("H126").Value = 2546.5833
$ws.Range("I126").Value = 2545.4644
$ws.Range("J126").Value = 2550.5
$ws.Range("K126").Value = 7636.3932
$ws.Range("L126").Value = 7651.5
$ws.Range("M126").Value = -5166.3932
$ws.Range("N126").Value = -12591.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4027.12
$ws.Range("I122").Value = 1194.8182
$ws.Range("J122").Value = 6252.5
$ws.Range("K122").Value = 3584.4546
$ws.Range("L122").Value = 18757.5
$ws.Range("M122").Value = -1134.4546
$ws.Range("N122").Value = -23657.5
$ws.Range("H132").Value = 6503.1
$ws.Range("I132").Value = 885.7646999999999
$ws.Range("J132").Value = 38334.668
$ws.Range("K132").Value = 2657.2941
$ws.Range("L132").Value = 115004.004
$ws.Range("M132").Value = -127.2941000000001
$ws.Range("N132").Value = -120064.004
$ws.Range("H135").Value = 46092.727
$ws.Range("J135").Value = 46092.727
$ws.Range("L135").Value = 46092.727
$ws.Range("N135").Value = -56232.727
$ws.Range("H140").Value = 43006.46
$ws.Range("J140").Value = 43006.46
$ws.Range("L140").Value = 43006.46
$ws.Range("N140").Value = -53366.46
$ws.Range("H141").Value = 47614.21
$ws.Range("J141").Value = 47614.21
$ws.Range("L141").Value = 47614.21
$ws.Range("N141").Value = -57974.21
